$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2401.4285
$ws.Range("I40").Value = 2900
$ws.Range("J40").Value = 2202
$ws.Range("K40").Value = 2900
$ws.Range("L40").Value = 2202
$ws.Range("M40").Value = -2725
$ws.Range("N40").Value = -2552
$ws.Range("H88").Value = 2068.6
$ws.Range("I88").Value = 1116.5
$ws.Range("J88").Value = 2703.3333
$ws.Range("K88").Value = 1116.5
$ws.Range("L88").Value = 2703.3333
$ws.Range("M88").Value = -710.5
$ws.Range("N88").Value = -3515.3333
$ws.Range("H91").Value = 2068.6
$ws.Range("I91").Value = 1116.5
$ws.Range("J91").Value = 2703.3333
$ws.Range("K91").Value = 1116.5
$ws.Range("L91").Value = 2703.3333
$ws.Range("M91").Value = 287.5
$ws.Range("N91").Value = -5511.3333
$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("L99").ClearContents()
$ws.Range("M99").ClearContents()
$ws.Range("N99").Value = 0
$ws.Range("H137").Value = 4857.6484
$ws.Range("I137").Value = 5103.75
$ws.Range("J137").Value = 4670.143
$ws.Range("K137").Value = 15311.25
$ws.Range("L137").Value = 14010.429
$ws.Range("M137").Value = -12761.25
$ws.Range("N137").Value = -19110.429
$ws.Range("H138").Value = 4212.6665
$ws.Range("I138").Value = 2438
$ws.Range("J138").Value = 5100
$ws.Range("K138").Value = 7314
$ws.Range("L138").Value = 15300
$ws.Range("M138").Value = -2174
$ws.Range("N138").Value = -25580

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 0
$ws.Range("K19").Value = 0
$ws.Range("M19").ClearContents()
$ws.Range("H21").Value = 9000
$ws.Range("I21").Value = 9000
$ws.Range("K21").Value = 9000
$ws.Range("M21").Value = -8626
$ws.Range("H23").Value = 29633.334
$ws.Range("I23").Value = 38000
$ws.Range("J23").Value = 25450
$ws.Range("K23").Value = 38000
$ws.Range("L23").Value = 25450
$ws.Range("M23").Value = -37741
$ws.Range("N23").Value = -25968
$ws.Range("H26").Value = 799.5
$ws.Range("I26").Value = 799.5
$ws.Range("J26").Value = 0
$ws.Range("K26").Value = 799.5
$ws.Range("L26").Value = 0
$ws.Range("M26").ClearContents()
$ws.Range("N26").Value = -469.5
$ws.Range("H30").Value = 0
$ws.Range("I30").Value = 0
$ws.Range("J30").Value = 0
$ws.Range("K30").Value = 0
$ws.Range("L30").ClearContents()
$ws.Range("M30").ClearContents()
$ws.Range("N30").Value = 0
$ws.Range("H74").Value = 5326
$ws.Range("I74").Value = 2057.7778
$ws.Range("K74").Value = 2057.7778
$ws.Range("M74").Value = -1183.7778
$ws.Range("H77").Value = 5326
$ws.Range("I77").Value = 2057.7778
$ws.Range("K77").Value = 10288.889
$ws.Range("M77").Value = -5920.888999999999
$ws.Range("H97").Value = 1276.25
$ws.Range("I97").Value = 1276.25
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 1276.25
$ws.Range("L97").Value = 0
$ws.Range("M97").ClearContents()
$ws.Range("N97").Value = -780.25
$ws.Range("H132").Value = 3430.8928
$ws.Range("I132").Value = 2993.3
$ws.Range("J132").Value = 4524.875
$ws.Range("K132").Value = 8979.900000000001
$ws.Range("L132").Value = 13574.625
$ws.Range("M132").Value = -6449.900000000001
$ws.Range("N132").Value = -18634.625

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H30").Value = 0
$ws.Range("I30").Value = 0
$ws.Range("K30").Value = 0
$ws.Range("M30").ClearContents()
$ws.Range("H33").Value = 0
$ws.Range("J33").Value = 0
$ws.Range("L33").ClearContents()
$ws.Range("N33").Value = 0
$ws.Range("H99").Value = 1378.762
$ws.Range("I99").Value = 1044.3529
$ws.Range("J99").Value = 2800
$ws.Range("K99").Value = 1044.3529
$ws.Range("L99").Value = 2800
$ws.Range("M99").Value = 453.6470999999999
$ws.Range("N99").Value = -5796
$ws.Range("H106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("L106").ClearContents()
$ws.Range("N106").Value = 0

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H11").Value = 175
$ws.Range("I11").Value = 175
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 175
$ws.Range("L11").Value = 0
$ws.Range("M11").ClearContents()
$ws.Range("N11").Value = -35
$ws.Range("H22").Value = 427.33334
$ws.Range("I22").Value = 390.5
$ws.Range("J22").Value = 501
$ws.Range("K22").Value = 390.5
$ws.Range("L22").Value = 501
$ws.Range("M22").Value = -40.5
$ws.Range("N22").Value = -1201
$ws.Range("H31").Value = 3997.578
$ws.Range("I31").Value = 5990.1
$ws.Range("J31").Value = 2403.56
$ws.Range("K31").Value = 5990.1
$ws.Range("L31").Value = 2403.56
$ws.Range("M31").Value = -5695.1
$ws.Range("N31").Value = -2993.56
$ws.Range("H34").Value = 3997.578
$ws.Range("I34").Value = 5990.1
$ws.Range("J34").Value = 2403.56
$ws.Range("K34").Value = 5990.1
$ws.Range("L34").Value = 2403.56
$ws.Range("M34").Value = -5788.1
$ws.Range("N34").Value = -2807.56
$ws.Range("H53").Value = 31496.334
$ws.Range("J53").Value = 31496.334
$ws.Range("L53").Value = 31496.334
$ws.Range("N53").Value = -32710.334
$ws.Range("H111").Value = 79800
$ws.Range("J111").Value = 79800
$ws.Range("L111").Value = 79800
$ws.Range("N111").Value = -87980
$ws.Range("H118").Value = 57450
$ws.Range("J118").Value = 57450
$ws.Range("L118").Value = 57450
$ws.Range("N118").Value = -60764
$ws.Range("H122").Value = 6876.7915
$ws.Range("I122").Value = 2510.1875
$ws.Range("K122").Value = 7530.5625
$ws.Range("M122").Value = -5080.5625
$ws.Range("H132").Value = 3937.0344
$ws.Range("I132").Value = 3389.9048
$ws.Range("J132").Value = 5373.25
$ws.Range("K132").Value = 10169.7144
$ws.Range("L132").Value = 16119.75
$ws.Range("M132").Value = -7639.714399999999
$ws.Range("N132").Value = -21179.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 1479.3103
$ws.Range("I129").Value = 1104.2858
$ws.Range("J129").Value = 1598.6364
$ws.Range("K129").Value = 3312.8574
$ws.Range("L129").Value = 4795.9092
$ws.Range("M129").Value = 1687.1426
$ws.Range("N129").Value = -14795.9092

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 15025002
$ws.Range("J18").Value = 50000
$ws.Range("L18").Value = 50000
$ws.Range("N18").Value = -50586
$ws.Range("H39").Value = 0
$ws.Range("J39").Value = 0
$ws.Range("L39").ClearContents()
$ws.Range("N39").Value = 0
$ws.Range("H132").Value = 3070.8125
$ws.Range("I132").Value = 1963.5
$ws.Range("K132").Value = 5890.5
$ws.Range("M132").Value = -3360.5
$ws.Range("H135").Value = 64880
$ws.Range("J135").Value = 64880
$ws.Range("L135").Value = 64880
$ws.Range("N135").Value = -75020

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1361.16
$ws.Range("J22").Value = 1875.0714
$ws.Range("L22").Value = 1875.0714
$ws.Range("N22").Value = -2465.0714
$ws.Range("H27").Value = 1361.16
$ws.Range("J27").Value = 1875.0714
$ws.Range("L27").Value = 1875.0714
$ws.Range("N27").Value = -2089.0714
$ws.Range("H45").Value = 16500
$ws.Range("I45").Value = 9500
$ws.Range("K45").Value = 9500
$ws.Range("M45").Value = -9093
$ws.Range("H55").Value = 363908.72
$ws.Range("I55").Value = 1000103.5
$ws.Range("K55").Value = 1000103.5
$ws.Range("M55").Value = -999930.5
$ws.Range("H81").Value = 38181
$ws.Range("J81").Value = 38181
$ws.Range("L81").Value = 38181
$ws.Range("N81").Value = -40177
$ws.Range("H82").Value = 2585.7144
$ws.Range("J82").Value = 2820
$ws.Range("L82").Value = 2820
$ws.Range("N82").Value = -3542
$ws.Range("H84").Value = 38181
$ws.Range("J84").Value = 38181
$ws.Range("L84").Value = 114543
$ws.Range("N84").Value = -124527
$ws.Range("H85").Value = 2585.7144
$ws.Range("J85").Value = 2820
$ws.Range("L85").Value = 2820
$ws.Range("N85").Value = -5316
$ws.Range("H141").Value = 79060.5
$ws.Range("J141").Value = 79060.5
$ws.Range("L141").Value = 79060.5
$ws.Range("N141").Value = -89420.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 14288696
$ws.Range("I81").Value = 2092.5
$ws.Range("J81").Value = 33337502
$ws.Range("K81").Value = 4185
$ws.Range("L81").Value = 66675004
$ws.Range("M81").Value = -3124
$ws.Range("N81").Value = -66677126
$ws.Range("H84").Value = 14288696
$ws.Range("I84").Value = 2092.5
$ws.Range("J84").Value = 33337502
$ws.Range("K84").Value = 20925
$ws.Range("L84").Value = 333375020
$ws.Range("M84").Value = -15621
$ws.Range("N84").Value = -333385628
$ws.Range("H110").Value = 48500
$ws.Range("J110").Value = 48500
$ws.Range("L110").Value = 48500
$ws.Range("N110").Value = -56680
$ws.Range("H116").Value = 58000
$ws.Range("J116").Value = 58000
$ws.Range("L116").Value = 58000
$ws.Range("N116").Value = -67178
